$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-07-16 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-17 Monday", 2) | Out-Null

# Update the arithmetic table (20 rows x 5 columns)
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "75-24="
$t.Cell(1,2).Range.Text = "79-60="
$t.Cell(1,3).Range.Text = "40-11="
$t.Cell(1,4).Range.Text = "20-17="
$t.Cell(1,5).Range.Text = "52-28="
$t.Cell(2,1).Range.Text = "69+15="
$t.Cell(2,2).Range.Text = "59-31="
$t.Cell(2,3).Range.Text = "49-17="
$t.Cell(2,4).Range.Text = "65-7="
$t.Cell(2,5).Range.Text = "35+21="
$t.Cell(3,1).Range.Text = "98-39="
$t.Cell(3,2).Range.Text = "9+70="
$t.Cell(3,3).Range.Text = "49+32="
$t.Cell(3,4).Range.Text = "25+45="
$t.Cell(3,5).Range.Text = "63-24="
$t.Cell(4,1).Range.Text = "25+57="
$t.Cell(4,2).Range.Text = "70-66="
$t.Cell(4,3).Range.Text = "55+19="
$t.Cell(4,4).Range.Text = "29+7="
$t.Cell(4,5).Range.Text = "43-34="
$t.Cell(5,1).Range.Text = "8+25="
$t.Cell(5,2).Range.Text = "24-16="
$t.Cell(5,3).Range.Text = "3+20="
$t.Cell(5,4).Range.Text = "98-45="
$t.Cell(5,5).Range.Text = "72+25="
$t.Cell(6,1).Range.Text = "93-62="
$t.Cell(6,2).Range.Text = "20+49="
$t.Cell(6,3).Range.Text = "66-57="
$t.Cell(6,4).Range.Text = "67+15="
$t.Cell(6,5).Range.Text = "63-52="
$t.Cell(7,1).Range.Text = "7+65="
$t.Cell(7,2).Range.Text = "6+76="
$t.Cell(7,3).Range.Text = "12+1="
$t.Cell(7,4).Range.Text = "52+38="
$t.Cell(7,5).Range.Text = "93-87="
$t.Cell(8,1).Range.Text = "74-44="
$t.Cell(8,2).Range.Text = "97-40="
$t.Cell(8,3).Range.Text = "33+8="
$t.Cell(8,4).Range.Text = "2+51="
$t.Cell(8,5).Range.Text = "19+49="
$t.Cell(9,1).Range.Text = "34+33="
$t.Cell(9,2).Range.Text = "32-26="
$t.Cell(9,3).Range.Text = "52+20="
$t.Cell(9,4).Range.Text = "99-36="
$t.Cell(9,5).Range.Text = "25+25="
$t.Cell(10,1).Range.Text = "88-24="
$t.Cell(10,2).Range.Text = "16+5="
$t.Cell(10,3).Range.Text = "75-74="
$t.Cell(10,4).Range.Text = "97-21="
$t.Cell(10,5).Range.Text = "60+12="
$t.Cell(11,1).Range.Text = "75-56="
$t.Cell(11,2).Range.Text = "66-28="
$t.Cell(11,3).Range.Text = "69-55="
$t.Cell(11,4).Range.Text = "1+40="
$t.Cell(11,5).Range.Text = "79-43="
$t.Cell(12,1).Range.Text = "44+47="
$t.Cell(12,2).Range.Text = "12+76="
$t.Cell(12,3).Range.Text = "25+47="
$t.Cell(12,4).Range.Text = "63+25="
$t.Cell(12,5).Range.Text = "74-33="
$t.Cell(13,1).Range.Text = "76+15="
$t.Cell(13,2).Range.Text = "49+39="
$t.Cell(13,3).Range.Text = "92-46="
$t.Cell(13,4).Range.Text = "31+20="
$t.Cell(13,5).Range.Text = "6-0="
$t.Cell(14,1).Range.Text = "26+10="
$t.Cell(14,2).Range.Text = "59-6="
$t.Cell(14,3).Range.Text = "37+0="
$t.Cell(14,4).Range.Text = "80-78="
$t.Cell(14,5).Range.Text = "4+21="
$t.Cell(15,1).Range.Text = "43-26="
$t.Cell(15,2).Range.Text = "40+29="
$t.Cell(15,3).Range.Text = "10-4="
$t.Cell(15,4).Range.Text = "70+1="
$t.Cell(15,5).Range.Text = "30-16="
$t.Cell(16,1).Range.Text = "69-5="
$t.Cell(16,2).Range.Text = "91-36="
$t.Cell(16,3).Range.Text = "41-40="
$t.Cell(16,4).Range.Text = "11+79="
$t.Cell(16,5).Range.Text = "1+31="
$t.Cell(17,1).Range.Text = "17+58="
$t.Cell(17,2).Range.Text = "18+15="
$t.Cell(17,3).Range.Text = "26-1="
$t.Cell(17,4).Range.Text = "13+29="
$t.Cell(17,5).Range.Text = "19-17="
$t.Cell(18,1).Range.Text = "3+3="
$t.Cell(18,2).Range.Text = "18+18="
$t.Cell(18,3).Range.Text = "63-4="
$t.Cell(18,4).Range.Text = "75-66="
$t.Cell(18,5).Range.Text = "16+75="
$t.Cell(19,1).Range.Text = "24+39="
$t.Cell(19,2).Range.Text = "81-16="
$t.Cell(19,3).Range.Text = "57-42="
$t.Cell(19,4).Range.Text = "42-12="
$t.Cell(19,5).Range.Text = "70-66="
$t.Cell(20,1).Range.Text = "57-50="
$t.Cell(20,2).Range.Text = "67-48="
$t.Cell(20,3).Range.Text = "96-63="
$t.Cell(20,4).Range.Text = "14+66="
$t.Cell(20,5).Range.Text = "98-28="

Write-Output "done"
